$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -and $val.Contains(",")) {
        $parts = $val.Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $n = $trimmed.Length
        $reversedParts = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversedParts += $trimmed[$i]
        }
        $newVal = $reversedParts -join ", "
        $cell.Value = $newVal
    }
}
